$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D88").Value = 0.7254492243564907
$ws.Range("D89").Value = 0.7215746373564907
$ws.Range("D90").Value = 0.5311946523564907
$ws.Range("D91").Value = 0.5539812373564907
$ws.Range("C92").Value = 0.1753415943564907
$ws.Range("C93").Value = 0.2651053283564908
$ws.Range("C94").Value = 0.08763596535649075
$ws.Range("C95").Value = 0.1003532183564907
$ws.Range("C96").Value = -0.02418658464350926
$ws.Range("C97").Value = 0.2001520573564908
